$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (fixing bug counts in allAdsController.js)
$ws.Range("C8").Value = 14
$ws.Range("C9").Value = 36
$ws.Range("C11").Value = 9

# Recalculate so the SUM formula in C51 updates to reflect new values
$excel.Calculate()

# Update the view: clear the frozen/top-left scroll position and move selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C9").Select()
